# [REVERT] David's merge of April 9th
#
# Reverts the "tech" sheet's per-platform (android/iOS) availability flags
# for rows 6-15 back to TRUE/TRUE, and restores the view state (active
# sheet/tab + scroll position + selection) to what it was before that
# merge.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("global_settings")
$ws2 = $wb.Worksheets.Item("tech")

# --- Data changes -----------------------------------------------------
# Columns F ([android]) and G ([iOS]) on the "tech" sheet, rows 6-15,
# go from FALSE back to TRUE.
foreach ($r in 6..15) {
    $ws2.Range("F$r").Value = $true
    $ws2.Range("G$r").Value = $true
}

# --- View/selection changes --------------------------------------------
# "tech" is no longer the selected tab; its scroll position/selection move.
$ws2.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 2
$ws2.Range("F21:H21").Select() | Out-Null

# "global_settings" becomes the selected/active tab again; its scroll
# position moves while its selection (F20) is unchanged.
$ws1.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("F20").Select() | Out-Null
